$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$c = $ws.Range("D2")
$c.NumberFormat = "@"
$c.Value = "57.830.06"
$c.Style = "Normal"
$ws.Range("E2").Value = "  -3.41%  "

$c = $ws.Range("D3")
$c.NumberFormat = "@"
$c.Value = "2.289.35"
$c.Style = "Normal"
$ws.Range("E3").Value = "  -3.78%  "

$ws.Range("E4").Value = "  +0.02%  "

$c = $ws.Range("D5")
$c.NumberFormat = "@"
$c.Value = "533.31"
$c.Style = "Normal"
$ws.Range("E5").Value = "  -4.12%  "

$c = $ws.Range("D6")
$c.NumberFormat = "@"
$c.Value = "130.39"
$c.Style = "Normal"
$ws.Range("E6").Value = "  -2.55%  "

$ws.Range("E7").Value = "  +0.06%  "

$ws.Range("E8").Value = "  -1.27%  "

$c = $ws.Range("D9")
$c.NumberFormat = "@"
$c.Value = "2.289.55"
$c.Style = "Normal"
$ws.Range("E9").Value = "  -3.73%  "

$c = $ws.Range("D10")
$c.NumberFormat = "@"
$c.Value = "0.0995"
$c.Style = "Normal"
$ws.Range("E10").Value = "  -5.68%  "

$c = $ws.Range("D11")
$c.NumberFormat = "@"
$c.Value = "5.41"
$c.Style = "Normal"
$ws.Range("E11").Value = "  -4.46%  "

$ws.Range("E12").Value = "  -0.61%  "

$ws.Range("E13").Value = "  -3.94%  "

$c = $ws.Range("D14")
$c.NumberFormat = "@"
$c.Value = "23.47"
$c.Style = "Normal"
$ws.Range("E14").Value = "  -4.06%  "

$c = $ws.Range("D15")
$c.NumberFormat = "@"
$c.Value = "2.698.58"
$c.Style = "Normal"
$ws.Range("E15").Value = "  -3.71%  "

$c = $ws.Range("D16")
$c.NumberFormat = "@"
$c.Value = "57.798.72"
$c.Style = "Normal"
$ws.Range("E16").Value = "  -3.35%  "

$ws.Range("E17").Value = "  -4.47%  "

$c = $ws.Range("D18")
$c.NumberFormat = "@"
$c.Value = "2.270.27"
$c.Style = "Normal"
$ws.Range("E18").Value = "  -4.34%  "

$c = $ws.Range("D19")
$c.NumberFormat = "@"
$c.Value = "10.49"
$c.Style = "Normal"
$ws.Range("E19").Value = "  -5.79%  "

$c = $ws.Range("D20")
$c.NumberFormat = "@"
$c.Value = "4.21"
$c.Style = "Normal"
$ws.Range("E20").Value = "  -6.08%  "

$c = $ws.Range("D21")
$c.NumberFormat = "@"
$c.Value = "312.98"
$c.Style = "Normal"
$ws.Range("E21").Value = "  -2.29%  "

$c = $ws.Range("D22")
$c.NumberFormat = "@"
$c.Value = "6.35"
$c.Style = "Normal"
$ws.Range("E22").Value = "  -4.78%  "

$c = $ws.Range("D23")
$c.NumberFormat = "@"
$c.Value = "1.00"
$c.Style = "Normal"
$ws.Range("E23").Value = "  -0.09%  "

$c = $ws.Range("D24")
$c.NumberFormat = "@"
$c.Value = "62.38"
$c.Style = "Normal"
$ws.Range("E24").Value = "  -2.67%  "

$c = $ws.Range("D25")
$c.NumberFormat = "@"
$c.Value = "0.165"
$c.Style = "Normal"
$ws.Range("E25").Value = "  -5.18%  "

$c = $ws.Range("D26")
$c.NumberFormat = "@"
$c.Value = "1.00"
$c.Style = "Normal"
$ws.Range("E26").Value = "  +0.13%  "

$c = $ws.Range("D27")
$c.NumberFormat = "@"
$c.Value = "8.02"
$c.Style = "Normal"
$ws.Range("E27").Value = "  -4.73%  "

$ws.Range("E28").Value = "  -6.89%  "

$c = $ws.Range("D29")
$c.NumberFormat = "@"
$c.Value = "170.27"
$c.Style = "Normal"
$ws.Range("E29").Value = "  +0.21%  "

$ws.Range("E30").Value = "  -5.67%  "

$c = $ws.Range("D31")
$c.NumberFormat = "@"
$c.Value = "0.0₃0715"
$c.Style = "Normal"
$ws.Range("E31").Value = "  -6.09%  "

$c = $ws.Range("D32")
$c.NumberFormat = "@"
$c.Value = "5.75"
$c.Style = "Normal"
$ws.Range("E32").Value = "  -5.45%  "

$c = $ws.Range("D33")
$c.NumberFormat = "@"
$c.Value = "1.04"
$c.Style = "Normal"
$ws.Range("E33").Value = "  -6.22%  "

$c = $ws.Range("D34")
$c.NumberFormat = "@"
$c.Value = "0.378"
$c.Style = "Normal"
$ws.Range("E34").Value = "  -5.43%  "

$ws.Range("E35").Value = "  +0.06%  "

$c = $ws.Range("D36")
$c.NumberFormat = "@"
$c.Value = "17.69"
$c.Style = "Normal"
$ws.Range("E36").Value = "  -2.67%  "

$ws.Range("E37").Value = "  -0.04%  "

$ws.Range("E38").Value = "  -7.42%  "

$ws.Range("E39").Value = "  -6.00%  "

$c = $ws.Range("D40")
$c.NumberFormat = "@"
$c.Value = "38.10"
$c.Style = "Normal"
$ws.Range("E40").Value = "  -1.40%  "

$c = $ws.Range("D41")
$c.NumberFormat = "@"
$c.Value = "1.48"
$c.Style = "Normal"
$ws.Range("E41").Value = "  -6.39%  "

$c = $ws.Range("D42")
$c.NumberFormat = "@"
$c.Value = "140.47"
$c.Style = "Normal"
$ws.Range("E42").Value = "  -3.22%  "

$c = $ws.Range("D43")
$c.NumberFormat = "@"
$c.Value = "288.01"
$c.Style = "Normal"
$ws.Range("E43").Value = "  -9.55%  "

$c = $ws.Range("D44")
$c.NumberFormat = "@"
$c.Value = "3.40"
$c.Style = "Normal"
$ws.Range("E44").Value = "  -3.77%  "

$c = $ws.Range("D45")
$c.NumberFormat = "@"
$c.Value = "0.0947"
$c.Style = "Normal"
$ws.Range("E45").Value = "  -2.46%  "

$ws.Range("E46").Value = "  -2.79%  "

$c = $ws.Range("D47")
$c.NumberFormat = "@"
$c.Value = "0.553"
$c.Style = "Normal"
$ws.Range("E47").Value = "  -3.33%  "

$c = $ws.Range("D48")
$c.NumberFormat = "@"
$c.Value = "18.08"
$c.Style = "Normal"
$ws.Range("E48").Value = "  -8.23%  "

$ws.Range("E49").Value = "  -3.69%  "

$c = $ws.Range("D50")
$c.NumberFormat = "@"
$c.Value = "10.93"
$c.Style = "Normal"
$ws.Range("E50").Value = "  -1.10%  "

$c = $ws.Range("D51")
$c.NumberFormat = "@"
$c.Value = "0.0₆0200"
$c.Style = "Normal"
$ws.Range("E51").Value = "  +83.98%  "
